$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing rows 2-6 down to 3-7
$ws.Rows(2).Insert()

# The inserted row inherits the header row formatting; strip it back to the default (no explicit style)
$ws.Range("A2:AO2").ClearFormats()

# Protect text-like columns (League/Date/Time/Home/Away) for the brand-new row 2 from Excel's
# automatic date/time literal-to-serial conversion by forcing Text format while we write them.
$ws.Range("A2:E2").NumberFormat = "@"

# Row 2: Japanese J League 3 - Kanazawa vs Matsumoto (07:00:00)
$ws.Cells.Item(2,1).Value = 'Japanese J League 3'
$ws.Cells.Item(2,2).Value = '2025-10-15'
$ws.Cells.Item(2,3).Value = '07:00:00'
$ws.Cells.Item(2,4).Value = 'Kanazawa'
$ws.Cells.Item(2,5).Value = 'Matsumoto'
$ws.Cells.Item(2,6).Value = 1.01
$ws.Cells.Item(2,7).Value = 1000
$ws.Cells.Item(2,8).Value = 1.01
$ws.Cells.Item(2,9).Value = 1000
$ws.Cells.Item(2,10).Value = 1.02
$ws.Cells.Item(2,11).Value = 1000
$ws.Cells.Item(2,12).Value = 1.01
$ws.Cells.Item(2,13).Value = 1.01
$ws.Cells.Item(2,14).Value = 1.28
$ws.Cells.Item(2,15).Value = 1.01
$ws.Cells.Item(2,16).Value = 1.28
$ws.Cells.Item(2,17).Value = 1.31
$ws.Cells.Item(2,18).Value = 1.18
$ws.Cells.Item(2,19).Value = 1.31
$ws.Cells.Item(2,20).Value = 1.01
$ws.Cells.Item(2,21).Value = 1.01
$ws.Cells.Item(2,22).Value = 1.01
$ws.Cells.Item(2,23).Value = 1.01
$ws.Cells.Item(2,24).Value = 1000
$ws.Cells.Item(2,25).Value = 1000
$ws.Cells.Item(2,26).Value = 1000
$ws.Cells.Item(2,27).Value = 1000
$ws.Cells.Item(2,28).Value = 1000
$ws.Cells.Item(2,29).Value = 1000
$ws.Cells.Item(2,30).Value = 1000
$ws.Cells.Item(2,31).Value = 1000
$ws.Cells.Item(2,32).Value = 1000
$ws.Cells.Item(2,33).Value = 1000
$ws.Cells.Item(2,34).Value = 1000
$ws.Cells.Item(2,35).Value = 1000
$ws.Cells.Item(2,36).Value = 1000
$ws.Cells.Item(2,37).Value = 1000
$ws.Cells.Item(2,38).Value = 1000
$ws.Cells.Item(2,39).Value = 1000
$ws.Cells.Item(2,40).Value = 1000
$ws.Cells.Item(2,41).Value = 1000

# Restore row 2 to default (General) formatting/style now that the literal text is committed
$ws.Range("A2:E2").ClearFormats()

# Row 3: refreshed odds (League/Date/Time/Home/Away already correct after the row shift)
$ws.Cells.Item(3,6).Value = 1.29
$ws.Cells.Item(3,7).Value = 1.43
$ws.Cells.Item(3,8).Value = 10
$ws.Cells.Item(3,9).Value = 1000
$ws.Cells.Item(3,10).Value = 4.7
$ws.Cells.Item(3,11).Value = 1000
$ws.Cells.Item(3,12).Value = 0
$ws.Cells.Item(3,13).Value = 0
$ws.Cells.Item(3,14).Value = 1.5
$ws.Cells.Item(3,15).Value = 1.28
$ws.Cells.Item(3,16).Value = 1.25
$ws.Cells.Item(3,17).Value = 1.28
$ws.Cells.Item(3,18).Value = 0
$ws.Cells.Item(3,19).Value = 0
$ws.Cells.Item(3,20).Value = 0
$ws.Cells.Item(3,21).Value = 0
$ws.Cells.Item(3,22).Value = 0
$ws.Cells.Item(3,23).Value = 0
$ws.Cells.Item(3,24).Value = 0
$ws.Cells.Item(3,25).Value = 0
$ws.Cells.Item(3,26).Value = 0
$ws.Cells.Item(3,27).Value = 0
$ws.Cells.Item(3,28).Value = 0
$ws.Cells.Item(3,29).Value = 0
$ws.Cells.Item(3,30).Value = 0
$ws.Cells.Item(3,31).Value = 0
$ws.Cells.Item(3,32).Value = 0
$ws.Cells.Item(3,33).Value = 0
$ws.Cells.Item(3,34).Value = 0
$ws.Cells.Item(3,35).Value = 0
$ws.Cells.Item(3,36).Value = 0
$ws.Cells.Item(3,37).Value = 0
$ws.Cells.Item(3,38).Value = 0
$ws.Cells.Item(3,39).Value = 0
$ws.Cells.Item(3,40).Value = 0
$ws.Cells.Item(3,41).Value = 0

# Row 4: refreshed odds (League/Date/Time/Home/Away already correct after the row shift)
$ws.Cells.Item(4,6).Value = 4
$ws.Cells.Item(4,7).Value = 1000
$ws.Cells.Item(4,8).Value = 2.04
$ws.Cells.Item(4,9).Value = 2.28
$ws.Cells.Item(4,10).Value = 3.15
$ws.Cells.Item(4,11).Value = 950
$ws.Cells.Item(4,12).Value = 0
$ws.Cells.Item(4,13).Value = 0
$ws.Cells.Item(4,14).Value = 0
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = 1.25
$ws.Cells.Item(4,17).Value = 1.01
$ws.Cells.Item(4,18).Value = 0
$ws.Cells.Item(4,19).Value = 0
$ws.Cells.Item(4,20).Value = 0
$ws.Cells.Item(4,21).Value = 0
$ws.Cells.Item(4,22).Value = 0
$ws.Cells.Item(4,23).Value = 0
$ws.Cells.Item(4,24).Value = 0
$ws.Cells.Item(4,25).Value = 0
$ws.Cells.Item(4,26).Value = 0
$ws.Cells.Item(4,27).Value = 0
$ws.Cells.Item(4,28).Value = 0
$ws.Cells.Item(4,29).Value = 0
$ws.Cells.Item(4,30).Value = 0
$ws.Cells.Item(4,31).Value = 0
$ws.Cells.Item(4,32).Value = 0
$ws.Cells.Item(4,33).Value = 0
$ws.Cells.Item(4,34).Value = 0
$ws.Cells.Item(4,35).Value = 0
$ws.Cells.Item(4,36).Value = 0
$ws.Cells.Item(4,37).Value = 0
$ws.Cells.Item(4,38).Value = 0
$ws.Cells.Item(4,39).Value = 0
$ws.Cells.Item(4,40).Value = 0
$ws.Cells.Item(4,41).Value = 0

# Row 5: refreshed odds (League/Date/Time/Home/Away already correct after the row shift)
$ws.Cells.Item(5,6).Value = 1.04
$ws.Cells.Item(5,7).Value = 980
$ws.Cells.Item(5,8).Value = 1.04
$ws.Cells.Item(5,9).Value = 1000
$ws.Cells.Item(5,10).Value = 1.01
$ws.Cells.Item(5,11).Value = 980
$ws.Cells.Item(5,12).Value = 0
$ws.Cells.Item(5,13).Value = 0
$ws.Cells.Item(5,14).Value = 0
$ws.Cells.Item(5,15).Value = 0
$ws.Cells.Item(5,16).Value = 1.25
$ws.Cells.Item(5,17).Value = 1.01
$ws.Cells.Item(5,18).Value = 0
$ws.Cells.Item(5,19).Value = 0
$ws.Cells.Item(5,20).Value = 0
$ws.Cells.Item(5,21).Value = 0
$ws.Cells.Item(5,22).Value = 0
$ws.Cells.Item(5,23).Value = 0
$ws.Cells.Item(5,24).Value = 0
$ws.Cells.Item(5,25).Value = 0
$ws.Cells.Item(5,26).Value = 0
$ws.Cells.Item(5,27).Value = 0
$ws.Cells.Item(5,28).Value = 0
$ws.Cells.Item(5,29).Value = 0
$ws.Cells.Item(5,30).Value = 0
$ws.Cells.Item(5,31).Value = 0
$ws.Cells.Item(5,32).Value = 0
$ws.Cells.Item(5,33).Value = 0
$ws.Cells.Item(5,34).Value = 0
$ws.Cells.Item(5,35).Value = 0
$ws.Cells.Item(5,36).Value = 0
$ws.Cells.Item(5,37).Value = 0
$ws.Cells.Item(5,38).Value = 0
$ws.Cells.Item(5,39).Value = 0
$ws.Cells.Item(5,40).Value = 0
$ws.Cells.Item(5,41).Value = 0

# Row 6: refreshed odds (League/Date/Time/Home/Away already correct after the row shift)
$ws.Cells.Item(6,6).Value = 1.09
$ws.Cells.Item(6,7).Value = 2.68
$ws.Cells.Item(6,8).Value = 1.42
$ws.Cells.Item(6,9).Value = 980
$ws.Cells.Item(6,10).Value = 2.96
$ws.Cells.Item(6,11).Value = 950
$ws.Cells.Item(6,12).Value = 0
$ws.Cells.Item(6,13).Value = 0
$ws.Cells.Item(6,14).Value = 0
$ws.Cells.Item(6,15).Value = 0
$ws.Cells.Item(6,16).Value = 1.25
$ws.Cells.Item(6,17).Value = 1.01
$ws.Cells.Item(6,18).Value = 0
$ws.Cells.Item(6,19).Value = 0
$ws.Cells.Item(6,20).Value = 0
$ws.Cells.Item(6,21).Value = 0
$ws.Cells.Item(6,22).Value = 0
$ws.Cells.Item(6,23).Value = 0
$ws.Cells.Item(6,24).Value = 0
$ws.Cells.Item(6,25).Value = 0
$ws.Cells.Item(6,26).Value = 0
$ws.Cells.Item(6,27).Value = 0
$ws.Cells.Item(6,28).Value = 0
$ws.Cells.Item(6,29).Value = 0
$ws.Cells.Item(6,30).Value = 0
$ws.Cells.Item(6,31).Value = 0
$ws.Cells.Item(6,32).Value = 0
$ws.Cells.Item(6,33).Value = 0
$ws.Cells.Item(6,34).Value = 0
$ws.Cells.Item(6,35).Value = 0
$ws.Cells.Item(6,36).Value = 0
$ws.Cells.Item(6,37).Value = 0
$ws.Cells.Item(6,38).Value = 0
$ws.Cells.Item(6,39).Value = 0
$ws.Cells.Item(6,40).Value = 0
$ws.Cells.Item(6,41).Value = 0

# Row 7: refreshed odds (League/Date/Time/Home/Away already correct after the row shift)
$ws.Cells.Item(7,6).Value = 1.04
$ws.Cells.Item(7,7).Value = 1000
$ws.Cells.Item(7,8).Value = 1.04
$ws.Cells.Item(7,9).Value = 1000
$ws.Cells.Item(7,10).Value = 1.01
$ws.Cells.Item(7,11).Value = 980
$ws.Cells.Item(7,12).Value = 0
$ws.Cells.Item(7,13).Value = 0
$ws.Cells.Item(7,14).Value = 0
$ws.Cells.Item(7,15).Value = 0
$ws.Cells.Item(7,16).Value = 1.24
$ws.Cells.Item(7,17).Value = 1.01
$ws.Cells.Item(7,18).Value = 0
$ws.Cells.Item(7,19).Value = 0
$ws.Cells.Item(7,20).Value = 0
$ws.Cells.Item(7,21).Value = 0
$ws.Cells.Item(7,22).Value = 0
$ws.Cells.Item(7,23).Value = 0
$ws.Cells.Item(7,24).Value = 0
$ws.Cells.Item(7,25).Value = 0
$ws.Cells.Item(7,26).Value = 0
$ws.Cells.Item(7,27).Value = 0
$ws.Cells.Item(7,28).Value = 0
$ws.Cells.Item(7,29).Value = 0
$ws.Cells.Item(7,30).Value = 0
$ws.Cells.Item(7,31).Value = 0
$ws.Cells.Item(7,32).Value = 0
$ws.Cells.Item(7,33).Value = 0
$ws.Cells.Item(7,34).Value = 0
$ws.Cells.Item(7,35).Value = 0
$ws.Cells.Item(7,36).Value = 0
$ws.Cells.Item(7,37).Value = 0
$ws.Cells.Item(7,38).Value = 0
$ws.Cells.Item(7,39).Value = 0
$ws.Cells.Item(7,40).Value = 0
$ws.Cells.Item(7,41).Value = 0

